$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to keep the literal text representation (matching the
    # original inlineStr cells) instead of letting Excel auto-coerce a
    # numeric-looking string into a real number. Style is restored to
    # "Normal" afterwards so no stray number-format style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "65.057.16"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.183.62"
$ws.Range("E3").Value = "  +1.33%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.25%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "615.63"
$ws.Range("E5").Value = "  +1.60%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "147.97"
$ws.Range("E6").Value = "  -1.38%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "3.182.12"
$ws.Range("E8").Value = "  +1.30%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.30%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.31%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "5.51"
$ws.Range("E11").Value = "  -1.55%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.13%  "

# Row 13 - ShibaInu
Set-TextValue $ws.Range("D13") "0.0000263"
$ws.Range("E13").Value = "  +1.12%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "36.10"
$ws.Range("E14").Value = "  -2.33%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.706.16"
$ws.Range("E15").Value = "  +1.40%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +3.13%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "65.041.97"
$ws.Range("E17").Value = "  +1.23%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "3.183.39"
$ws.Range("E18").Value = "  +1.43%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -0.85%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "484.15"
$ws.Range("E20").Value = "  +0.18%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("D21") "14.79"
$ws.Range("E21").Value = "  +1.02%  "

# Row 22 - Polygon
Set-TextValue $ws.Range("D22") "0.727"
$ws.Range("E22").Value = "  +1.89%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +2.57%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D24") "13.99"
$ws.Range("E24").Value = "  +1.05%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "84.69"
$ws.Range("E25").Value = "  +0.67%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.02%  "

# Row 27 - now PancakeSwap (was RenderToken)
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D27") "2.84"
$ws.Range("E27").Value = "  -3.38%  "

# Row 28 - now RenderToken (was PancakeSwap)
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D28") "8.73"
$ws.Range("E28").Value = "  +1.76%  "

# Row 29 - NEARProtocol
Set-TextValue $ws.Range("D29") "7.13"
$ws.Range("E29").Value = "  +0.52%  "

# Row 30 - Hedera
Set-TextValue $ws.Range("D30") "0.121"
$ws.Range("E30").Value = "  -4.29%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -5.16%  "

# Row 32 - Stacks
Set-TextValue $ws.Range("D32") "2.73"
$ws.Range("E32").Value = "  -0.22%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.36%  "

# Row 34 - EthereumClassic
Set-TextValue $ws.Range("D34") "26.85"
$ws.Range("E34").Value = "  +0.28%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  +2.66%  "

# Row 36 - PEPE
Set-TextValue $ws.Range("D36") "0.0₃0797"
$ws.Range("E36").Value = "  +5.38%  "

# Row 37 - Filecoin
Set-TextValue $ws.Range("D37") "6.06"
$ws.Range("E37").Value = "  -0.79%  "

# Row 38 - dogwifhat
Set-TextValue $ws.Range("D38") "3.20"
$ws.Range("E38").Value = "  -1.86%  "

# Row 39 - OKB
Set-TextValue $ws.Range("D39") "53.25"
$ws.Range("E39").Value = "  -2.18%  "

# Row 40 - Bittensor
Set-TextValue $ws.Range("D40") "470.64"
$ws.Range("E40").Value = "  +3.61%  "

# Row 41 - VeChain
Set-TextValue $ws.Range("D41") "0.0403"
$ws.Range("E41").Value = "  +0.32%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -2.62%  "

# Row 43 - Cosmos
$ws.Range("E43").Value = "  -0.94%  "

# Row 44 - Maker
Set-TextValue $ws.Range("D44") "2.869.46"
$ws.Range("E44").Value = "  -0.69%  "

# Row 45 - Fetch.AI
Set-TextValue $ws.Range("D45") "2.35"
$ws.Range("E45").Value = "  +0.99%  "

# Row 46 - TheGraph
Set-TextValue $ws.Range("D46") "0.271"
$ws.Range("E46").Value = "  -0.79%  "

# Row 47 - ThetaToken
Set-TextValue $ws.Range("D47") "2.48"
$ws.Range("E47").Value = "  +6.39%  "

# Row 48 - Arweave
Set-TextValue $ws.Range("D48") "37.44"
$ws.Range("E48").Value = "  +10.71%  "

# Row 49 - InjectiveProtocol
Set-TextValue $ws.Range("D49") "26.98"
$ws.Range("E49").Value = "  +0.76%  "

# Row 50 - USDe
$ws.Range("E50").Value = "  +0.06%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -0.50%  "
